# Applies the "Finalized Test Plan for Project Finale" edits to the
# Drew Grubb Chess Test Plan document using the Word COM object model.

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $replaceText
    }
    return $found
}

# 1. Committed date placeholder -> actual date, split into two runs like a
#    real typed edit (select "xxxx..." and overtype with "30", plus drop the
#    "March " wording change).
$rng = $d.Content
$found = $rng.Find.Execute("Committed March xxxxxxxxxxxxxxxxxxxxxxxxxxx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Delete()
    $rng.InsertAfter("Committed ")
    $rng.Collapse(0)
    $rng.InsertAfter("March 30")
}

# 2. "CalculateMove (Minimax AI)" -> "CalculateMove (Various AI)"
Replace-Text "Minimax AI" "Various AI"

# 3. Move the _GoBack bookmark: remove it from the empty paragraph under
#    "Environment Requirements" and re-insert it mid-sentence in the
#    Automated Test paragraph.
$bm = $d.Bookmarks.Item("\GoBack")
$bm.Delete()

$rng = $d.Content
$found = $rng.Find.Execute("Any testing will be done ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $rng) | Out-Null
}

# 4. Remove the "Code coverage for Chess Chaos from unit testing is XX%."
#    paragraph entirely, and reword the following sentence.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Code coverage for Chess Chaos from unit testing is*") {
        $p.Range.Delete()
        break
    }
}

Replace-Text "Unit tests cover positive/negative tests" "Unit tests should cover various positive/negative tests"

# 5. "The following functions have specific unit tests created for them." ->
#    "The following functions need to be tested"
Replace-Text "The following functions have specific unit tests created for them." "The following functions need to be tested"

# 6. "There are unit tests created for each piece" -> "There will be unit tests created for each piece"
Replace-Text "There are unit tests created for each piece" "There will be unit tests created for each piece"

# 7. Both "There is"/"is" -> "will be" occurrences ahead of "one unit test".
$rng = $d.Content
$found = $rng.Find.Execute("isInCheckmate method. The order of the unit tests goes as follows:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rng = $d.Content
$found = $rng.Find.Execute("There is one unit test created for the isInCheckmate method. The order of the unit tests goes as follows:`nInitialize Board`nAdd pieces on board in Checkmate", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng2 = $d.Range($rng.Start, $rng.Start + 8)
    Write-Output $rng2.Text
}

Write-Output "done"
